# Update vm_pu.xlsx results for "case with 380 kV done"
# Applies new bus voltage magnitude (vm_pu) values for rows 2-25
# across columns B-N (column H stays empty, column G stays 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$newValues = @(
    @(1.02, 1.028597543587176, 1.031838044503245, 1.038066542135152, 1.048042148412875, 1.030172257777299, 1.033748977520374, 1.034645007092213, 1.040855634039331, 1.050803098583792, 1.035217019641119),
    @(1.02, 1.029622203118716, 1.032583262993766, 1.039050544069596, 1.049279246135172, 1.030314164999929, 1.034413928410534, 1.035199042699446, 1.041649112735315, 1.051851061320019, 1.035882914837813),
    @(1.02, 1.030285380811341, 1.033065294662082, 1.039687812707354, 1.050080761641904, 1.030404537247507, 1.034843784612176, 1.035556703654162, 1.042162486939805, 1.052529631417169, 1.036313381484566),
    @(1.02, 1.030564217666142, 1.033267898143713, 1.039955852665668, 1.05041796568621, 1.030442181966128, 1.035024397306142, 1.03570686335123, 1.042378294971867, 1.052815014405787, 1.036494250669274),
    @(1.02, 1.030611037767298, 1.03330191363509, 1.040000865498851, 1.050474598214999, 1.030448482276188, 1.035054717165469, 1.035732064024003, 1.042414529234319, 1.052862938041578, 1.036524613586279),
    @(1.02, 1.030289106497868, 1.033068002028309, 1.039691393751234, 1.050085266408332, 1.030405041625036, 1.03484619835669, 1.035558710885758, 1.042165370634481, 1.052533444278012, 1.036315798656875),
    @(1.02, 1.028943800076057, 1.032089930208016, 1.038398975387088, 1.048460018120048, 1.030220516260514, 1.033973785782783, 1.034832419018298, 1.041123806080217, 1.051157165595441, 1.035442147157045),
    @(1.02, 1.026574381924955, 1.030365143219506, 1.036125822130882, 1.045604006226156, 1.029884259785076, 1.03243334665154, 1.033546212861271, 1.039287989515238, 1.048735569530847, 1.033899520425545),
    @(1.02, 1.024995563605313, 1.029214460042988, 1.034613261507427, 1.0437052821901, 1.029652646381341, 1.031404296372009, 1.032684479036042, 1.038063822062729, 1.047123574720999, 1.032869008776583),
    @(1.02, 1.024312104717664, 1.028716015504728, 1.033958990397401, 1.042884362426254, 1.029550593979043, 1.030958213518186, 1.032310333722511, 1.037533676843126, 1.046426130735561, 1.032422292433942),
    @(1.02, 1.024058264446956, 1.028530843012427, 1.033716066886819, 1.042579622115611, 1.029512422690712, 1.030792443950632, 1.032171208341735, 1.037336746425301, 1.046167152880554, 1.032256287454588),
    @(1.02, 1.02411271283803, 1.028570564396064, 1.033768170191051, 1.042644981544926, 1.029520622525534, 1.030828005455529, 1.032201058041413, 1.037378989160136, 1.046222700717928, 1.0322918994609),
    @(1.02, 1.024291121656495, 1.028700709655751, 1.033938908184661, 1.042859168708861, 1.029547444122186, 1.030944512473451, 1.032298836656173, 1.037517398735072, 1.046404721835566, 1.032408571932152),
    @(1.02, 1.024401048761015, 1.028780892757051, 1.034044119005002, 1.042991161100626, 1.02956393475499, 1.031016286421313, 1.032359061271213, 1.037602676001815, 1.046516882161588, 1.032480447807253),
    @(1.02, 1.025040926290474, 1.029247536199832, 1.034656697611942, 1.043759790019894, 1.029659382157982, 1.031433890959752, 1.032709288608236, 1.038099004527503, 1.047169873562767, 1.032898645392035),
    @(1.02, 1.025442352549954, 1.029540198639631, 1.035041133481146, 1.044242262339266, 1.029718782179401, 1.031695710240739, 1.032928707155109, 1.038410318876231, 1.047579627640846, 1.033160836486432),
    @(1.02, 1.025676515065086, 1.029710885200116, 1.03526543379296, 1.044523799981807, 1.029753259154709, 1.031848377004227, 1.033056592910737, 1.038591896258299, 1.047818684352938, 1.033313720054236),
    @(1.02, 1.0257563613847, 1.029769081748348, 1.035341925488107, 1.044619817360917, 1.029764986060033, 1.031900424279355, 1.033100182120014, 1.038653808263087, 1.047900205697892, 1.033365841242465),
    @(1.02, 1.025399281516301, 1.029508800636739, 1.034999880387893, 1.044190485247625, 1.029712426706515, 1.03166762447254, 1.032905175689967, 1.038376918517654, 1.04753565931331, 1.033132710833222),
    @(1.02, 1.024238583968456, 1.02866238588768, 1.033888627286971, 1.042796090801395, 1.029539553129247, 1.030910206116596, 1.032270047460373, 1.037476640865043, 1.046351118839658, 1.032374216856335),
    @(1.02, 1.02350896224541, 1.02813004970895, 1.033190528036749, 1.041920454646745, 1.029429330575362, 1.030433555993556, 1.031869842819863, 1.036910537650978, 1.045606836473815, 1.031896889835448),
    @(1.02, 1.023895733920107, 1.028412266271419, 1.033560547727172, 1.042384544209046, 1.029487906576911, 1.030686278049734, 1.032082081606751, 1.03721064562837, 1.046001348741881, 1.032149970785933),
    @(1.02, 1.025418743408586, 1.029522988099864, 1.035018520682909, 1.044213880716518, 1.029715298996572, 1.031680315370462, 1.032915808846986, 1.038392010723539, 1.0475555265406, 1.033145419753674),
    @(1.02, 1.027186793260871, 1.030811190983846, 1.036712981146791, 1.046341421252767, 1.029972503389357, 1.032831956942925, 1.033879481659104, 1.039762643951763, 1.049361186237906, 1.034298696789266)
)

$startRow = 2

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range($columns[$j] + $rowIndex).Value = $rowValues[$j]
    }
}
